$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
# ($wb.ActiveSheet also resolves to Sheet1 here, but look it up explicitly by
# name so the script stays correct even if the active sheet selection differs.)

# Insert a new blank row above row 6 (shifts existing rows 6-18 down to 7-19)
$ws.Rows("6:6").Insert()

# Update G5 (Scenario 1 / S001 employer) value from "Sams" to new value "Taco bell"
$ws.Range("G5").Value = "Taco bell"

# Update the selection to match the target state
$ws.Range("K6").Select()
